$d = $word.ActiveDocument

# 0a. Title paragraph: remove the hanging indent, and append "Jim Baldwin" on a
#     manual line break inside the same (bold) run as the title text.
$p1 = $d.Paragraphs(1)
$p1.LeftIndent = 0
$p1.RightIndent = 0
$p1.FirstLineIndent = -0.0001
$titleText = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$titleText.InsertAfter([char]11 + "Jim Baldwin")

# 0b. Second paragraph (intro "The ... support site ..." paragraph): pick up a
#     matching (zeroed) indent definition that it did not have before.
$p2 = $d.Paragraphs(2)
$p2.LeftIndent = 0
$p2.RightIndent = 0
$p2.FirstLineIndent = -0.0001

# 1. Merge "Access to online Knowledge " + "B" + "ase" -> "Access to online Knowledge Base"
$r = $d.Content
$r.Find.Execute("Access to online Knowledge Base", $true, $false, $false, $false, $false, $true, 1, $false, "Access to online Knowledge Base", 2)

# 2. Merge "Knowledgebase – Owned by " + "Director of Professional Services"
$r = $d.Content
$s2 = "Knowledgebase – Owned by Director of Professional Services"
$r.Find.Execute($s2, $true, $false, $false, $false, $false, $true, 1, $false, $s2, 2)

# 3. Merge "End of Life Policy – static text page – " + "Dir. Sales, Dir. Operations"
$r = $d.Content
$s3 = "End of Life Policy – static text page – Dir. Sales, Dir. Operations"
$r.Find.Execute($s3, $true, $false, $false, $false, $false, $true, 1, $false, $s3, 2)

# 4. Merge "Partner Portal – Owned by " + "Channel Marketing Manager"
$r = $d.Content
$s4 = "Partner Portal – Owned by Channel Marketing Manager"
$r.Find.Execute($s4, $true, $false, $false, $false, $false, $true, 1, $false, $s4, 2)

# 5. Merge "Sends an email..." split across 3 runs
$r = $d.Content
$sendsText = "Sends an email to [support email address]. This could be a form, but it needs to send an email to [support email address], because the CRM system uses that to open a support ticket automatically."
$r.Find.Execute($sendsText, $true, $false, $false, $false, $false, $true, 1, $false, $sendsText, 2)

# 6. Merge "Sales/Marketing Portal – owned by " + "Director of Sales, Director of Marketing"
$r = $d.Content
$s6 = "Sales/Marketing Portal – owned by Director of Sales, Director of Marketing"
$r.Find.Execute($s6, $true, $false, $false, $false, $false, $true, 1, $false, $s6, 2)
